$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 106
$ws.Range("F3").Value = 1259
$ws.Range("F5").Value = 932
$ws.Range("F6").Value = 1682
$ws.Range("F7").Value = 363
$ws.Range("F8").Value = 1127
$ws.Range("F9").Value = 45
$ws.Range("F10").Value = 99
$ws.Range("F11").Value = 251
$ws.Range("F12").Value = 15
$ws.Range("F13").Value = 75
$ws.Range("F14").Value = 614
$ws.Range("F15").Value = 119
$ws.Range("F17").Value = 21
$ws.Range("F19").Value = 313
$ws.Range("F20").Value = 77
$ws.Range("F21").Value = 631
$ws.Range("F22").Value = 618
$ws.Range("F23").Value = 106
$ws.Range("F24").Value = 26
$ws.Range("F27").Value = 56
$ws.Range("F29").Value = 237

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 8
$ws.Range("G5").Value = 118
$ws.Range("F8").Value = 81
$ws.Range("F11").Value = 18

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 106
$ws.Range("F4").Value = 1259
$ws.Range("F6").Value = 932
$ws.Range("F7").Value = 1682
$ws.Range("F8").Value = 363
$ws.Range("F9").Value = 1127
$ws.Range("F10").Value = 45
$ws.Range("F12").Value = 99
$ws.Range("F13").Value = 251
$ws.Range("F14").Value = 15
$ws.Range("F15").Value = 75
$ws.Range("F16").Value = 614
$ws.Range("F17").Value = 119
$ws.Range("F20").Value = 21
$ws.Range("F23").Value = 8
$ws.Range("G23").Value = 118
$ws.Range("F24").Value = 313
$ws.Range("F28").Value = 77
$ws.Range("F29").Value = 631
$ws.Range("F30").Value = 618
$ws.Range("F31").Value = 106
$ws.Range("F32").Value = 26
$ws.Range("F35").Value = 81
$ws.Range("F36").Value = 56
$ws.Range("F38").Value = 237
$ws.Range("F44").Value = 18
